$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7440765500068665
$ws.Range("B1").Value = 3.091638565063477
$ws.Range("C1").Value = 3.028201818466187
$ws.Range("D1").Value = 2.410148859024048
$ws.Range("E1").Value = 1.492595791816711
